# Update crypto price/volume data in the worksheet to reflect the latest
# scrape from coinranking.com (GitHub Actions scheduled update).
#
# Note: Price values are stored as plain text (not numbers) in this sheet
# (e.g. "24.548.76", "0.00001313", "5.830") so every write to column D
# forces the cell's number format to Text ("@") first. Without this, Excel's
# COM layer would auto-coerce numeric-looking strings into real numbers and
# silently mangle them (drop significant trailing zeros, switch to
# scientific notation, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 45 and 46 swapped their coin identity (NEARProtocol and Decentraland
# traded ranking positions) along with their price/volume data.
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.556"
$ws.Range("E45").Value = "  -5.97%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7079"
$ws.Range("E46").Value = "  -4.73%  "

# Price / volume(1h) refresh for all other rows.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.548.76"
$ws.Range("E2").Value = "  -1.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.671.06"
$ws.Range("E3").Value = "  -2.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.87"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3954"
$ws.Range("E7").Value = "  -1.66%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3939"
$ws.Range("E8").Value = "  -2.65%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.395"
$ws.Range("E10").Value = "  -5.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.41"
$ws.Range("E11").Value = "  -6.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08636"
$ws.Range("E12").Value = "  -1.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.36"
$ws.Range("E13").Value = "  -3.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.297"
$ws.Range("E14").Value = "  -2.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001313"
$ws.Range("E15").Value = "  -2.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.675"
$ws.Range("E16").Value = "  -4.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.679.62"
$ws.Range("E17").Value = "  +3.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.93"
$ws.Range("E18").Value = "  -1.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07011"
$ws.Range("E19").Value = "  -2.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "21.18"
$ws.Range("E20").Value = "  +1.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.073"
$ws.Range("E21").Value = "  -2.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.87"
$ws.Range("E23").Value = "  -4.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.561.61"
$ws.Range("E24").Value = "  -1.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.343"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.758"
$ws.Range("E26").Value = "  -4.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.99"
$ws.Range("E27").Value = "  -0.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.830"
$ws.Range("E28").Value = "  -9.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "159.11"
$ws.Range("E29").Value = "  -1.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "145.52"
$ws.Range("E30").Value = "  +1.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.269"
$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.527"
$ws.Range("E32").Value = "  +10.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.861.86"
$ws.Range("E33").Value = "  -1.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.03077"
$ws.Range("E34").Value = "  -3.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08245"
$ws.Range("E35").Value = "  -5.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.899"
$ws.Range("E36").Value = "  -4.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2802"
$ws.Range("E37").Value = "  -2.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9896"
$ws.Range("E38").Value = "  -3.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09632"
$ws.Range("E39").Value = "  +1.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.516"
$ws.Range("E40").Value = "  +2.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.26"
$ws.Range("E41").Value = "  -5.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7853"
$ws.Range("E42").Value = "  -6.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.51"
$ws.Range("E43").Value = "  -4.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.50"
$ws.Range("E44").Value = "  -6.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.174"
$ws.Range("E47").Value = "  -1.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08634"
$ws.Range("E48").Value = "  +2.74%  "

$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.324"
$ws.Range("E50").Value = "  -3.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.70"
$ws.Range("E51").Value = "  -2.06%  "
